# =========================================================================
# Craniosynostosis.xlsx — "Refined metadata to be additional tab"
#   1. Re-stamp the "time_taken" column (F2:F83) on the "data" sheet with
#      the refreshed query timestamps.
#   2. Add a new "metadata" sheet (placed right after "data") carrying the
#      panel-level query metadata that used to live only implicitly.
# =========================================================================

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# -------------------------------------------------------------------
# 1) Refresh the "time_taken" timestamps (F2:F83) on the "data" sheet
# -------------------------------------------------------------------
$ws1.Range("F2").Value = "2021-10-05 14:33:38.909141"
$ws1.Range("F3").Value = "2021-10-05 14:33:38.909149"
$ws1.Range("F4").Value = "2021-10-05 14:33:38.909152"
$ws1.Range("F5").Value = "2021-10-05 14:33:38.909155"
$ws1.Range("F6").Value = "2021-10-05 14:33:38.909158"
$ws1.Range("F7").Value = "2021-10-05 14:33:38.909160"
$ws1.Range("F8").Value = "2021-10-05 14:33:38.909163"
$ws1.Range("F9").Value = "2021-10-05 14:33:38.909165"
$ws1.Range("F10").Value = "2021-10-05 14:33:38.909168"
$ws1.Range("F11").Value = "2021-10-05 14:33:38.909170"
$ws1.Range("F12").Value = "2021-10-05 14:33:38.909173"
$ws1.Range("F13").Value = "2021-10-05 14:33:38.909175"
$ws1.Range("F14").Value = "2021-10-05 14:33:38.909178"
$ws1.Range("F15").Value = "2021-10-05 14:33:38.909180"
$ws1.Range("F16").Value = "2021-10-05 14:33:38.909183"
$ws1.Range("F17").Value = "2021-10-05 14:33:38.909185"
$ws1.Range("F18").Value = "2021-10-05 14:33:38.909188"
$ws1.Range("F19").Value = "2021-10-05 14:33:38.909191"
$ws1.Range("F20").Value = "2021-10-05 14:33:38.909193"
$ws1.Range("F21").Value = "2021-10-05 14:33:38.909195"
$ws1.Range("F22").Value = "2021-10-05 14:33:38.909198"
$ws1.Range("F23").Value = "2021-10-05 14:33:38.909200"
$ws1.Range("F24").Value = "2021-10-05 14:33:38.909203"
$ws1.Range("F25").Value = "2021-10-05 14:33:38.909205"
$ws1.Range("F26").Value = "2021-10-05 14:33:38.909208"
$ws1.Range("F27").Value = "2021-10-05 14:33:38.909211"
$ws1.Range("F28").Value = "2021-10-05 14:33:38.909213"
$ws1.Range("F29").Value = "2021-10-05 14:33:38.909216"
$ws1.Range("F30").Value = "2021-10-05 14:33:38.909218"
$ws1.Range("F31").Value = "2021-10-05 14:33:38.909221"
$ws1.Range("F32").Value = "2021-10-05 14:33:38.909223"
$ws1.Range("F33").Value = "2021-10-05 14:33:38.909225"
$ws1.Range("F34").Value = "2021-10-05 14:33:38.909228"
$ws1.Range("F35").Value = "2021-10-05 14:33:38.909231"
$ws1.Range("F36").Value = "2021-10-05 14:33:38.909233"
$ws1.Range("F37").Value = "2021-10-05 14:33:38.909235"
$ws1.Range("F38").Value = "2021-10-05 14:33:38.909238"
$ws1.Range("F39").Value = "2021-10-05 14:33:38.909240"
$ws1.Range("F40").Value = "2021-10-05 14:33:38.909243"
$ws1.Range("F41").Value = "2021-10-05 14:33:38.909245"
$ws1.Range("F42").Value = "2021-10-05 14:33:38.909248"
$ws1.Range("F43").Value = "2021-10-05 14:33:38.909250"
$ws1.Range("F44").Value = "2021-10-05 14:33:38.909253"
$ws1.Range("F45").Value = "2021-10-05 14:33:38.909255"
$ws1.Range("F46").Value = "2021-10-05 14:33:38.909258"
$ws1.Range("F47").Value = "2021-10-05 14:33:38.909260"
$ws1.Range("F48").Value = "2021-10-05 14:33:38.909262"
$ws1.Range("F49").Value = "2021-10-05 14:33:38.909265"
$ws1.Range("F50").Value = "2021-10-05 14:33:38.909267"
$ws1.Range("F51").Value = "2021-10-05 14:33:38.909270"
$ws1.Range("F52").Value = "2021-10-05 14:33:38.909272"
$ws1.Range("F53").Value = "2021-10-05 14:33:38.909274"
$ws1.Range("F54").Value = "2021-10-05 14:33:38.909277"
$ws1.Range("F55").Value = "2021-10-05 14:33:38.909280"
$ws1.Range("F56").Value = "2021-10-05 14:33:38.909282"
$ws1.Range("F57").Value = "2021-10-05 14:33:38.909285"
$ws1.Range("F58").Value = "2021-10-05 14:33:38.909287"
$ws1.Range("F59").Value = "2021-10-05 14:33:38.909290"
$ws1.Range("F60").Value = "2021-10-05 14:33:38.909292"
$ws1.Range("F61").Value = "2021-10-05 14:33:38.909295"
$ws1.Range("F62").Value = "2021-10-05 14:33:38.909297"
$ws1.Range("F63").Value = "2021-10-05 14:33:38.909300"
$ws1.Range("F64").Value = "2021-10-05 14:33:38.909302"
$ws1.Range("F65").Value = "2021-10-05 14:33:38.909305"
$ws1.Range("F66").Value = "2021-10-05 14:33:38.909308"
$ws1.Range("F67").Value = "2021-10-05 14:33:38.909311"
$ws1.Range("F68").Value = "2021-10-05 14:33:38.909314"
$ws1.Range("F69").Value = "2021-10-05 14:33:38.909316"
$ws1.Range("F70").Value = "2021-10-05 14:33:38.909319"
$ws1.Range("F71").Value = "2021-10-05 14:33:38.909321"
$ws1.Range("F72").Value = "2021-10-05 14:33:38.909323"
$ws1.Range("F73").Value = "2021-10-05 14:33:38.909326"
$ws1.Range("F74").Value = "2021-10-05 14:33:38.909328"
$ws1.Range("F75").Value = "2021-10-05 14:33:38.909331"
$ws1.Range("F76").Value = "2021-10-05 14:33:38.909333"
$ws1.Range("F77").Value = "2021-10-05 14:33:38.909336"
$ws1.Range("F78").Value = "2021-10-05 14:33:38.909341"
$ws1.Range("F79").Value = "2021-10-05 14:33:38.909344"
$ws1.Range("F80").Value = "2021-10-05 14:33:38.909346"
$ws1.Range("F81").Value = "2021-10-05 14:33:38.909349"
$ws1.Range("F82").Value = "2021-10-05 14:33:38.909351"
$ws1.Range("F83").Value = "2021-10-05 14:33:38.909354"

# -------------------------------------------------------------------
# 2) Add the new "metadata" sheet right after "data"
# -------------------------------------------------------------------
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "metadata"
$ws2.Cells.Clear()

# Header row (row 1)
$ws2.Range("B1").Value = "data_name"
$ws2.Range("C1").Value = "data_id"
$ws2.Range("D1").Value = "data_version"
$ws2.Range("E1").Value = "data_version_created"
$ws2.Range("F1").Value = "panel_query_time"
$ws2.Range("G1").Value = "panel_get_request"

# Reuse the "data" sheet's header formatting (bold / thin border / centered+top aligned)
$ws1.Range("B1:F1").Copy()
$ws2.Range("B1:F1").PasteSpecial(-4122)
$ws1.Range("F1").Copy()
$ws2.Range("G1").PasteSpecial(-4122)

# Data row (row 2)
$ws2.Range("A2").Value = 0
$ws1.Range("A2").Copy()
$ws2.Range("A2").PasteSpecial(-4122)

$ws2.Range("B2").Value = "Craniosynostosis"
$ws2.Range("C2").Value = 93
$ws2.Range("D2").Value = "'1.25"
$ws2.Range("E2").Value = "2021-09-05T02:33:36.669790Z"
$ws2.Range("F2").Value = "2021-10-05 14:33:38.905569"
$ws2.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/93/?format=json"

# Keep "data" as the active/selected sheet, matching the original workbook view
$ws1.Activate()
$null = $ws1.Range("A1").Select()
